$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '19.948.76'
$ws.Range("E2").Value = '  -6.98%  '

# Row 3
$ws.Range("D3").Value = '1.408.42'
$ws.Range("E3").Value = '  -7.99%  '

# Row 4
$ws.Range("D4").Value = '''1.002'
$ws.Range("E4").Value = '  +0.00%  '

# Row 5
$ws.Range("D5").Value = '''1.002'
$ws.Range("E5").Value = '  -0.03%  '

# Row 6
$ws.Range("D6").Value = '''275.85'

# Row 7
$ws.Range("E7").Value = '  -6.17%  '

# Row 8
$ws.Range("D8").Value = '''0.3106'
$ws.Range("E8").Value = '  -2.17%  '

# Row 9
$ws.Range("D9").Value = '''39.79'
$ws.Range("E9").Value = '  -6.58%  '

# Row 10
$ws.Range("D10").Value = '''1.029'
$ws.Range("E10").Value = '  -3.61%  '

# Row 11
$ws.Range("D11").Value = '''0.06493'
$ws.Range("E11").Value = '  -9.03%  '

# Row 12
$ws.Range("D12").Value = '''1.003'
$ws.Range("E12").Value = '  -0.03%  '

# Row 13
$ws.Range("D13").Value = '''5.501'
$ws.Range("E13").Value = '  -3.72%  '

# Row 14
$ws.Range("D14").Value = '''17.62'
$ws.Range("E14").Value = '  -2.91%  '

# Row 15
$ws.Range("D15").Value = '''6.182'
$ws.Range("E15").Value = '  -4.72%  '

# Row 16
$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").Value = '1.411.13'
$ws.Range("E16").Value = '  -7.84%  '

# Row 17
$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").Value = '''0.00001020'
$ws.Range("E17").Value = '  -6.00%  '

# Row 18
$ws.Range("D18").Value = '''0.05679'
$ws.Range("E18").Value = '  -13.91%  '

# Row 19
$ws.Range("D19").Value = '''1.001'
$ws.Range("E19").Value = '  -0.02%  '

# Row 20
$ws.Range("D20").Value = '''70.75'
$ws.Range("E20").Value = '  -14.57%  '

# Row 21
$ws.Range("D21").Value = '''5.611'
$ws.Range("E21").Value = '  -7.84%  '

# Row 22
$ws.Range("D22").Value = '''14.72'
$ws.Range("E22").Value = '  -4.27%  '

# Row 23
$ws.Range("D23").Value = '''10.93'
$ws.Range("E23").Value = '  +0.94%  '

# Row 24
$ws.Range("D24").Value = '''2.254'
$ws.Range("E24").Value = '  -5.00%  '

# Row 25
$ws.Range("D25").Value = '19.969.80'
$ws.Range("E25").Value = '  -6.93%  '

# Row 26
$ws.Range("D26").Value = '''2.261'
$ws.Range("E26").Value = '  -4.33%  '

# Row 27
$ws.Range("D27").Value = '''132.75'
$ws.Range("E27").Value = '  -10.84%  '

# Row 28
$ws.Range("D28").Value = '''17.20'
$ws.Range("E28").Value = '  -5.86%  '

# Row 29
$ws.Range("D29").Value = '1.570.11'

# Row 30
$ws.Range("D30").Value = '''109.57'
$ws.Range("E30").Value = '  -5.91%  '

# Row 31
$ws.Range("D31").Value = '''3.971'
$ws.Range("E31").Value = '  -17.64%  '

# Row 32
$ws.Range("D32").Value = '''5.282'
$ws.Range("E32").Value = '  -12.41%  '

# Row 33
$ws.Range("D33").Value = '''0.8175'
$ws.Range("E33").Value = '  -14.22%  '

# Row 34
$ws.Range("D34").Value = '''0.07677'
$ws.Range("E34").Value = '  -3.88%  '

# Row 35
$ws.Range("D35").Value = '''8.360'
$ws.Range("E35").Value = '  -1.00%  '

# Row 36
$ws.Range("D36").Value = '''1.480'
$ws.Range("E36").Value = '  -0.77%  '

# Row 37
$ws.Range("B37").Value = 'InternetComputer(DFINITY)'
$ws.Range("C37").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D37").Value = '''4.921'
$ws.Range("E37").Value = '  -4.20%  '

# Row 38
$ws.Range("B38").Value = 'Hedera'
$ws.Range("C38").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D38").Value = '''0.05881'
$ws.Range("E38").Value = '  -0.33%  '

# Row 39
$ws.Range("D39").Value = '''1.001'
$ws.Range("E39").Value = '  -0.10%  '

# Row 40
$ws.Range("D40").Value = '''0.02074'
$ws.Range("E40").Value = '  -5.10%  '

# Row 41
$ws.Range("D41").Value = '''10.46'
$ws.Range("E41").Value = '  -6.89%  '

# Row 42
$ws.Range("D42").Value = '''0.1900'
$ws.Range("E42").Value = '  -5.65%  '

# Row 43
$ws.Range("D43").Value = '''1.094'
$ws.Range("E43").Value = '  -6.88%  '

# Row 44
$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").Value = '''12.37'
$ws.Range("E44").Value = '  -5.19%  '

# Row 45
$ws.Range("B45").Value = 'TheSandbox'
$ws.Range("C45").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D45").Value = '''0.5306'
$ws.Range("E45").Value = '  -7.38%  '

# Row 46
$ws.Range("D46").Value = '''3.534'
$ws.Range("E46").Value = '  -4.69%  '

# Row 47
$ws.Range("D47").Value = '''0.5178'
$ws.Range("E47").Value = '  -6.55%  '

# Row 48
$ws.Range("D48").Value = '''114.61'
$ws.Range("E48").Value = '  -0.75%  '

# Row 49
$ws.Range("D49").Value = '''1.768'
$ws.Range("E49").Value = '  -6.49%  '

# Row 50
$ws.Range("D50").Value = '''1.035'
$ws.Range("E50").Value = '  -10.61%  '

# Row 51
$ws.Range("E51").Value = '  -0.03%  '
